$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows rows 3-6 (D, L, M, N, O, P, Q, R, S, T columns) being
# cyclically rotated among each other:
#   new row 3 <- old row 5
#   new row 4 <- old row 6
#   new row 5 <- old row 4
#   new row 6 <- old row 3
# Capture the "old" values first, then write them to their new rows.

$rows = 3, 4, 5, 6
$cols = "D", "L", "M", "N", "O", "P", "Q", "R", "S", "T"

$old = @{}
foreach ($r in $rows) {
    $old[$r] = @{}
    foreach ($c in $cols) {
        $old[$r][$c] = $ws.Range("$c$r").Value2
    }
}

$mapping = @{ 3 = 5; 4 = 6; 5 = 4; 6 = 3 }

foreach ($newRow in $rows) {
    $srcRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $old[$srcRow][$c]
    }
}
